$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

# Force text storage for the date/time columns so Excel doesn't
# auto-convert "2025-10-16" / "15:22:56" into date/time serials.
# Apply the text number format before writing the value, then restore
# the cell's style back to Normal (General) so no residual formatting
# is left behind on the new row.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-16"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "15:22:56"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,679.5098"
